$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.790.70'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.287.62'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '115.52'
$ws.Range('E5').Value = '  +17.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.41'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.620'
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.69'
$ws.Range('E10').Value = '  +8.41%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.04'
$ws.Range('E12').Value = '  +14.18%  '
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.82'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Value = '2.634.20'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.875'
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').Value = '2.282.13'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '43.661.80'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000109'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.99'
$ws.Range('E20').Value = '  +12.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.60'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.14'
$ws.Range('E23').Value = '  +11.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '233.41'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.96'
$ws.Range('E25').Value = '  +3.65%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.71'
$ws.Range('E27').Value = '  +4.58%  '
$ws.Range('E28').Value = '  +57.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '42.12'
$ws.Range('E29').Value = '  +9.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.38'
$ws.Range('E30').Value = '  -1.99%  '
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '173.72'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0935'
$ws.Range('E33').Value = '  +4.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.58'
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('E35').Value = '  +5.55%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.77'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0361'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.86'
$ws.Range('E40').Value = '  +9.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.52'
$ws.Range('E41').Value = '  +19.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '75.04'
$ws.Range('E42').Value = '  +15.76%  '
$ws.Range('E43').Value = '  +4.41%  '
$ws.Range('E44').Value = '  +2.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.38'
$ws.Range('E45').Value = '  +21.94%  '
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.40'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.71'
$ws.Range('E49').Value = '  +4.36%  '
$ws.Range('E50').Value = '  +3.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1000'
$ws.Range('E51').Value = '  -2.10%  '
